$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt")
$ws.Activate()

# --- Row 52: "UI For In-Game Stats" task moved from "this week" to "complete", hours 2 -> 6
$ws.Range("C52").Value = 6
$ws.Range("D52").Value = "complete"

# --- Row 53: "In-Game Menu" moved from "this week" to "complete" (hours already 7)
$ws.Range("D53").Value = "complete"

# --- Row 57: "Expand On Start Screen" hours spent 2 -> 7
$ws.Range("C57").Value = 7

# --- Insert a new row before row 58 for a "Training" task line
$ws.Range("A58").EntireRow.Insert()

# Populate the newly inserted row 58
$ws.Range("A58").Value = "Training"
$ws.Range("B58").Value = 6
$ws.Range("C58").Value = 6
$ws.Range("D58").Value = "complete"

# The old row 58 (now row 59, "Dr-BC Mode") keeps its values; nothing to change there.

# The "totals" row (now row 60) gained a value in its D column
$ws.Range("D60").Value = "planned"

# --- Best effort: fix up the conditional formatting range that pointed at the
# single cell D58 (old "Dr-BC Mode" row) so it keeps following that row (now D59)
$fc = $ws.Range("D58").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("D59"))

# --- Restore a sensible view state (scrolled down a bit further, selection on D61)
$ws.Range("A34").Select()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D61").Select()

$excel.Calculate()
